# Atualização de bases das ligas, do dia: 16-06-2024 às 07:16
# Swap the per-match data (columns B through AD) between pairs of rows
# while leaving column A (the sequential row index) untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chile Primera Division")

# Pairs of worksheet rows whose B:AD contents need to be swapped.
$pairs = @(
    @(102, 103),
    @(105, 106),
    @(108, 110),
    @(115, 116),
    @(118, 122),
    @(119, 121),
    @(155, 156)
)

$cols = @("B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC","AD")

foreach ($pair in $pairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    foreach ($col in $cols) {
        $cell1 = $ws.Range("$col$r1")
        $cell2 = $ws.Range("$col$r2")

        $value1 = $cell1.Value2
        $value2 = $cell2.Value2

        $cell1.Value2 = $value2
        $cell2.Value2 = $value1
    }
}
